$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the units reported in row 4 (the "Example" row) from kW-hr to W-hr.
# F4 holds the unit label; G4, J4, N4 and O4 hold quantities expressed in
# kW-hr (or derived from kW-hr) that need to be rescaled by 1/1000 to match
# the new W-hr unit label.
$ws.Range("F4").Value = "W-hr"
$ws.Range("G4").Value = 0.0036
$ws.Range("J4").Value = 0.0001355
$ws.Range("N4").Value = 8.83333333333333
$ws.Range("O4").Value = 81.8

# Update the selected cell on the sheet (cosmetic change captured in the diff).
$ws.Range("L10").Select()
